$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = @{ AB = 5.39; AC = 4.88 }
    3 = @{ AB = 4.62; AC = 4.21 }
    4 = @{ AB = 5.61; AC = 2.5 }
    5 = @{ AB = 2.94; AC = 3.35 }
    6 = @{ AB = 1.11; AC = 0.57 }
    7 = @{ AB = 8.050000000000001; AC = 9.199999999999999 }
    8 = @{ AB = 0.62; AC = 0.82 }
    9 = @{ AB = 3.57; AC = 4.94 }
}

foreach ($row in $values.Keys) {
    $ws.Range("AB$row").Value = $values[$row].AB
    $ws.Range("AC$row").Value = $values[$row].AC
}
